$d = $word.ActiveDocument

function Insert-LineBreakAfter($paraIndex, $anchorText) {
    # Locate $anchorText inside the given paragraph, collapse the found
    # range to its end, then insert a soft line break (w:br) there by
    # assigning the Word "manual line break" character (vertical tab,
    # chr 11) to the (now zero-length) range's Text property.
    $p = $d.Paragraphs($paraIndex)
    $r = $p.Range
    $found = $r.Find.Execute($anchorText, $true, $false, $false, $false, $false, `
                              $true, 1, $false, "", 0)
    if (-not $found) {
        Write-Output "NOT FOUND: [$anchorText] in paragraph $paraIndex"
        return
    }
    $r.Collapse(0)
    $r.Text = [char]11
}

# --- "Programa resumido" (English / italic) paragraph -----------------
# Split "...6. Production scales" + "Case studies (biotechnological ...)"
Insert-LineBreakAfter 12 "6. Production scales"

# --- "Programa" (English / italic) paragraph ---------------------------
# Split into the 9 numbered items, processed last-anchor-first so that
# earlier anchor text in the paragraph is never disturbed by a prior
# insertion.
$anchors15 = @(
    "otechnological processes.",
    "for Engineering in Brazil",
    " the Biochemical Engineer",
    "t recovery, among others.",
    "es, main unit operations.",
    "atory, pilot, industrial.",
    "technological processes).",
    "nnovation in Engineering."
)
for ($i = $anchors15.Length - 1; $i -ge 0; $i--) {
    Insert-LineBreakAfter 15 $anchors15[$i]
}

# --- "Bibliografia" paragraph -------------------------------------------
$anchors19 = @(
    "a Industrial – Engenharia",
    "lo: Edgard Blücher, 2001.",
    "epts. Second edition. New",
    "Jersey: PrenticeHall,",
    "2002."
)
for ($i = $anchors19.Length - 1; $i -ge 0; $i--) {
    Insert-LineBreakAfter 19 $anchors19[$i]
}
